# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.060.26"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3
$ws.Range("D3").Value = "2.395.13"
$ws.Range("E3").Value = "  +6.46%  "

# Row 4
$ws.Range("E4").Value = "  -0.47%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.12"
$ws.Range("E5").Value = "  +11.98%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.15"
$ws.Range("E6").Value = "  -6.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.649"
$ws.Range("E7").Value = "  +3.17%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.654"
$ws.Range("E9").Value = "  +7.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.05"
$ws.Range("E10").Value = "  -4.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0939"
$ws.Range("E11").Value = "  +1.29%  "

# Row 12
$ws.Range("E12").Value = "  -2.47%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.06"
$ws.Range("E13").Value = "  -1.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.23"
$ws.Range("E14").Value = "  +13.11%  "

# Row 15
$ws.Range("E15").Value = "  +1.43%  "

# Row 16
$ws.Range("D16").Value = "2.752.90"
$ws.Range("E16").Value = "  +6.38%  "

# Row 17
$ws.Range("D17").Value = "2.391.13"
$ws.Range("E17").Value = "  +2.48%  "

# Row 18
$ws.Range("D18").Value = "43.151.06"
$ws.Range("E18").Value = "  +1.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  +7.80%  "

# Row 20
$ws.Range("E20").Value = "  +2.00%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.57"
$ws.Range("E21").Value = "  +2.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.75"
$ws.Range("E22").Value = "  +7.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "272.59"
$ws.Range("E23").Value = "  +7.47%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  -1.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.64"
$ws.Range("E25").Value = "  +7.95%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.79"
$ws.Range("E26").Value = "  +2.10%  "

# Row 27
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.02"
$ws.Range("E28").Value = "  +3.31%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.44"
$ws.Range("E29").Value = "  -0.04%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("E30").Value = "  -1.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.73"
$ws.Range("E31").Value = "  -0.63%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.16"
$ws.Range("E32").Value = "  +0.16%  "

# Row 33
$ws.Range("E33").Value = "  +4.53%  "

# Row 34
$ws.Range("E34").Value = "  +4.08%  "

# Row 35
$ws.Range("E35").Value = "  +5.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.95"
$ws.Range("E36").Value = "  -2.21%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.16"
$ws.Range("E37").Value = "  -1.63%  "

# Row 38
$ws.Range("E38").Value = "  -2.60%  "

# Row 39
$ws.Range("E39").Value = "  +3.93%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +16.25%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("E41").Value = "  +19.59%  "

# Row 42
$ws.Range("E42").Value = "  +1.05%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.87"
$ws.Range("E43").Value = "  -2.97%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.93"
$ws.Range("E44").Value = "  +14.29%  "

# Row 45
$ws.Range("E45").Value = "  +0.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.41"
$ws.Range("E46").Value = "  -0.62%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.95"
$ws.Range("E47").Value = "  +49.23%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.36"
$ws.Range("E48").Value = "  +8.23%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.55"
$ws.Range("E49").Value = "  +0.61%  "

# Row 50
$ws.Range("E50").Value = "  +1.83%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.495"
$ws.Range("E51").Value = "  +12.78%  "
